# Kenntnisse.xlsx - update skill rating and cursor position
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 = "R" skill, column C = quantitative Beurteilung value: 2 -> 2.5
$ws.Range("C6").Value = 2.5

# Update the last active selection cell to C18
$ws.Range("C18").Select()
